$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: PartOf
$ws.Range("A10").Value = "PartOf"

# Update CapableOf template wording (row 3)
$ws.Range("B3").Value = "Its special ability is to <CapableOf>"

# New row 11: MadeOf
$ws.Range("A11").Value = "MadeOf"

# New row 12: ReceivesAction
$ws.Range("A12").Value = "ReceivesAction"

# Update UsedFor template wording (row 6)
$ws.Range("B6").Value = "It is used for <UsedFor>"

# Fill in templates for the three new rows
$ws.Range("B10").Value = "It is part of <PartOf>"
$ws.Range("B11").Value = "It is made of <MadeOf>"
$ws.Range("B12").Value = "It can be <ReceivesAction>"

# Update selection to match final state
$ws.Range("C12").Select()
